$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(12, 9).Value = 'sv'
$ws.Cells.Item(12, 10).Value = 'Statement-opinion'
$ws.Cells.Item(27, 9).Value = 'b'
$ws.Cells.Item(27, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(48, 9).Value = 'aa'
$ws.Cells.Item(48, 10).Value = 'Agree/Accept'
$ws.Cells.Item(54, 9).Value = 'sd'
$ws.Cells.Item(54, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(79, 9).Value = 'ba'
$ws.Cells.Item(79, 10).Value = 'Appreciation'
$ws.Cells.Item(80, 9).Value = 'sv'
$ws.Cells.Item(80, 10).Value = 'Statement-opinion'
$ws.Cells.Item(98, 9).Value = 'sv'
$ws.Cells.Item(98, 10).Value = 'Statement-opinion'
$ws.Cells.Item(99, 9).Value = 'b'
$ws.Cells.Item(99, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(102, 9).Value = 'b'
$ws.Cells.Item(102, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(108, 9).Value = 'sv'
$ws.Cells.Item(108, 10).Value = 'Statement-opinion'
$ws.Cells.Item(153, 9).Value = 'sv'
$ws.Cells.Item(153, 10).Value = 'Statement-opinion'
$ws.Cells.Item(155, 9).Value = 'b'
$ws.Cells.Item(155, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(176, 9).Value = 'sd'
$ws.Cells.Item(176, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(177, 9).Value = 'b'
$ws.Cells.Item(177, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(187, 9).Value = 'b'
$ws.Cells.Item(187, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(210, 9).Value = 'sv'
$ws.Cells.Item(210, 10).Value = 'Statement-opinion'
$ws.Cells.Item(225, 9).Value = 'sd'
$ws.Cells.Item(225, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(232, 9).Value = 'sv'
$ws.Cells.Item(232, 10).Value = 'Statement-opinion'
$ws.Cells.Item(235, 9).Value = 'sd'
$ws.Cells.Item(235, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(238, 9).Value = 'sd'
$ws.Cells.Item(238, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(253, 9).Value = 'sv'
$ws.Cells.Item(253, 10).Value = 'Statement-opinion'
$ws.Cells.Item(258, 9).Value = 'b'
$ws.Cells.Item(258, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(275, 9).Value = 'sv'
$ws.Cells.Item(275, 10).Value = 'Statement-opinion'
$ws.Cells.Item(279, 9).Value = 'b'
$ws.Cells.Item(279, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(283, 9).Value = 'sd'
$ws.Cells.Item(283, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(284, 9).Value = 'sv'
$ws.Cells.Item(284, 10).Value = 'Statement-opinion'
$ws.Cells.Item(304, 9).Value = 'sv'
$ws.Cells.Item(304, 10).Value = 'Statement-opinion'
$ws.Cells.Item(309, 9).Value = 'sd'
$ws.Cells.Item(309, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(341, 9).Value = 'ba'
$ws.Cells.Item(341, 10).Value = 'Appreciation'
$ws.Cells.Item(348, 9).Value = 'sd'
$ws.Cells.Item(348, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(352, 9).Value = 'ba'
$ws.Cells.Item(352, 10).Value = 'Appreciation'
$ws.Cells.Item(356, 9).Value = 'b'
$ws.Cells.Item(356, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(385, 9).Value = 'sv'
$ws.Cells.Item(385, 10).Value = 'Statement-opinion'
$ws.Cells.Item(386, 9).Value = 'b'
$ws.Cells.Item(386, 10).Value = 'Acknowledge (Backchannel)'
